$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue 'D2' '27.396.72'
Set-TextValue 'E2' '  +1.26%  '
Set-TextValue 'D3' '1.862.37'
Set-TextValue 'E3' '  +1.78%  '
Set-TextValue 'E4' '  -0.49%  '
Set-TextValue 'D5' '315.74'
Set-TextValue 'E5' '  +1.04%  '
Set-TextValue 'E6' '  -0.37%  '
Set-TextValue 'D7' '0.4619'
Set-TextValue 'E8' '  +0.55%  '
Set-TextValue 'D9' '0.07321'
Set-TextValue 'E9' '  -0.31%  '
Set-TextValue 'D10' '0.8901'
Set-TextValue 'E10' '  +2.13%  '
Set-TextValue 'D11' '20.06'
Set-TextValue 'E11' '  +1.28%  '
Set-TextValue 'D12' '0.07858'
Set-TextValue 'E12' '  -0.84%  '
Set-TextValue 'D13' '1.803.13'
Set-TextValue 'E13' '  +1.55%  '
Set-TextValue 'E14' '  +1.03%  '
Set-TextValue 'D15' '6.560'
Set-TextValue 'E15' '  -0.11%  '
Set-TextValue 'D16' '91.88'
Set-TextValue 'E16' '  +0.11%  '
Set-TextValue 'D17' '1.005'
Set-TextValue 'E17' '  -0.42%  '
Set-TextValue 'D18' '0.000008969'
Set-TextValue 'E18' '  +1.05%  '
Set-TextValue 'E19' '  -0.41%  '
Set-TextValue 'D20' '14.82'
Set-TextValue 'E20' '  +0.65%  '
Set-TextValue 'D21' '27.414.31'
Set-TextValue 'E21' '  +0.98%  '
Set-TextValue 'D22' '5.138'
Set-TextValue 'E22' '  +0.28%  '
Set-TextValue 'D23' '10.58'
Set-TextValue 'D24' '2.056.81'
Set-TextValue 'E24' '  -0.71%  '
Set-TextValue 'D25' '1.938'
Set-TextValue 'E25' '  +5.22%  '
Set-TextValue 'D26' '152.37'
Set-TextValue 'E26' '  -0.12%  '
Set-TextValue 'D27' '18.48'
Set-TextValue 'E27' '  -0.09%  '
Set-TextValue 'D28' '2.053'
Set-TextValue 'E28' '  -0.96%  '
Set-TextValue 'B29' 'BitcoinCash'
Set-TextValue 'C29' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D29' '116.43'
Set-TextValue 'E29' '  +1.09%  '
Set-TextValue 'B30' 'InternetComputer(DFINITY)'
Set-TextValue 'C30' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D30' '5.099'
Set-TextValue 'E30' '  -0.53%  '
Set-TextValue 'D31' '0.08849'
Set-TextValue 'E31' '  -0.26%  '
Set-TextValue 'D32' '3.091'
Set-TextValue 'E32' '  +3.89%  '
Set-TextValue 'D33' '0.7686'
Set-TextValue 'E33' '  +4.88%  '
Set-TextValue 'E34' '  +3.68%  '
Set-TextValue 'D35' '4.524'
Set-TextValue 'E35' '  +1.95%  '
Set-TextValue 'E36' '  +10.56%  '
Set-TextValue 'D37' '1.082'
Set-TextValue 'E37' '  +0.76%  '
Set-TextValue 'D38' '0.01964'
Set-TextValue 'E38' '  +1.47%  '
Set-TextValue 'B39' 'Hedera'
Set-TextValue 'C39' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D39' '0.05256'
Set-TextValue 'E39' '  +0.29%  '
Set-TextValue 'B40' 'MXToken'
Set-TextValue 'C40' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D40' '2.986'
Set-TextValue 'E40' '  +1.24%  '
Set-TextValue 'D41' '7.088'
Set-TextValue 'E41' '  -1.27%  '
Set-TextValue 'D42' '0.5153'
Set-TextValue 'E42' '  -0.05%  '
Set-TextValue 'D43' '0.1646'
Set-TextValue 'E43' '  +0.95%  '
Set-TextValue 'D44' '8.419'
Set-TextValue 'E44' '  +2.09%  '
Set-TextValue 'D45' '0.4810'
Set-TextValue 'E45' '  -0.47%  '
Set-TextValue 'D46' '10.42'
Set-TextValue 'E46' '  +1.63%  '
Set-TextValue 'D48' '103.40'
Set-TextValue 'E48' '  +0.96%  '
Set-TextValue 'D49' '1.649'
Set-TextValue 'E49' '  +1.36%  '
Set-TextValue 'D50' '0.06229'
Set-TextValue 'E50' '  +0.16%  '
Set-TextValue 'D51' '65.95'
Set-TextValue 'E51' '  +2.06%  '
